# Generate Report for Handback
#
# - Overview sheet: status text "Ready for handoff" -> "Handed back: in sync with en-US"
#   (achieved implicitly by updating the shared "Status" text on zh-cn/de-de sheets,
#   since Overview!E2/F2/E3/F3 mirror the same status string)
# - zh-cn / de-de sheets: fill in "Latest Target File" (I) / "Latest Handback File" (J) /
#   "Latest Handback DateTime" (K) for each of the two source files, now that the
#   localized files have been handed back.

$wb = $excel.ActiveWorkbook

$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8d6bba4a04a7e9a6510cd5270951575cec72098a/e2e/"

$files = @(
    @{ Name = "325f1c3d-cfab-42a6-bd50-5afd46afa46a"; Hash = "c036d781f6998aa0c646b5222fe723dd52e86868" },
    @{ Name = "99110e90-c253-4cc6-b28b-2cf5c814571b"; Hash = "c102cb304ac75afde20fddd1e01e7ac9464df026" }
)

$langSheets = @(
    @{ SheetName = "zh-cn"; Lang = "zh-cn"; HandbackTime = "2016-08-22 06:59:45" },
    @{ SheetName = "de-de"; Lang = "de-de"; HandbackTime = "2016-08-22 06:59:52" }
)

foreach ($langInfo in $langSheets) {
    $ws = $wb.Worksheets.Item($langInfo.SheetName)

    $row = 2
    foreach ($f in $files) {
        $mdName = "$($f.Name).md"
        $xlfName = "$($f.Name).$($f.Hash).$($langInfo.Lang).xlf"
        $mdUrl = "$baseUrl$mdName"

        # Status column (C) is now "handed back, in sync with source"
        $ws.Range("C$row").Value = "Handed back: in sync with en-US"

        # Latest Target File (I) -> hyperlink to the source .md, same as column A
        $ws.Range("I$row").Value = $mdName
        $ws.Hyperlinks.Add($ws.Range("I$row"), $mdUrl, [System.Type]::Missing, [System.Type]::Missing, $mdName)

        # Latest Handback File (J) -> localized xliff file name
        $ws.Range("J$row").Value = $xlfName

        # Latest Handback DateTime (K)
        $ws.Range("K$row").Value = $langInfo.HandbackTime

        $row = $row + 1
    }

    # Widen columns that now hold the longer "Handed back..." text / long file names
    $ws.Columns.Item(3).ColumnWidth = 29.144371396019366
    $ws.Columns.Item(9).ColumnWidth = 39.166666666666664
    $ws.Columns.Item(10).ColumnWidth = 39.166666666666664
}

# Overview sheet columns E/F also hold the status text and need the same widening
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 29.144371396019366
$overview.Columns.Item(6).ColumnWidth = 29.144371396019366
